$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 5 (the old "4.c" row) and column H (the old "4." column) -
# the correlation table loses the 4th variable entirely.
$ws.Rows.Item(5).Delete()
$ws.Columns.Item(8).Delete()

# Helper-free: force M/SD/ICC columns (B:D) to be stored as text (they were
# numeric before, but the new data is textual, e.g. ".49") while avoiding
# leaving a lingering number-format style on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 (previously "1.w") becomes "1.a"
$ws.Range("A2").Value = "1.a"
Set-TextValue $ws.Range("B2") "3.03"
Set-TextValue $ws.Range("C2") "0.99"
Set-TextValue $ws.Range("D2") ".49"
$ws.Range("E2").Value = "-"
$ws.Range("F2").Value = ".32**"
$ws.Range("G2").Value = ".39***"

# Row 3 (previously "2.a") becomes "2.b"
$ws.Range("A3").Value = "2.b"
Set-TextValue $ws.Range("B3") "1.95"
Set-TextValue $ws.Range("C3") "1.07"
Set-TextValue $ws.Range("D3") ".55"
$ws.Range("E3").Value = ".22***"
$ws.Range("F3").Value = "-"
$ws.Range("G3").Value = ".29**"

# Row 4 (previously "3.b") becomes "3.c"
$ws.Range("A4").Value = "3.c"
Set-TextValue $ws.Range("B4") "1.99"
Set-TextValue $ws.Range("C4") "1.06"
Set-TextValue $ws.Range("D4") ".47"
$ws.Range("E4").Value = ".25***"
$ws.Range("F4").Value = ".29***"
$ws.Range("G4").Value = "-"
